# feat: add 2022-Q1 data
#
# - Insert a new "2022-Q1" sheet (positioned right before the "总计" sheet)
#   holding the per-fund holdings detail for the new quarter. It is built
#   by duplicating the "总计" sheet (so it inherits the exact same header /
#   index-column cell style already used throughout the workbook) and then
#   overwriting its contents.
# - Update the "总计" (totals) sheet by adding a new leading row that
#   summarizes the 2022-Q1 quarter, shifting the previously-existing rows
#   down by one.

$wb = $excel.ActiveWorkbook

# Helper: write a value into a cell while forcing it to be stored as TEXT
# (matches the source data, which keeps numeric-looking strings like
# "3.16" / "98.58" as text rather than numbers), then drop the number
# format back to General so no stray formatting sticks around.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
}

# ---------------------------------------------------------------------
# 1. Duplicate "总计" to become the new "2022-Q1" sheet, positioned right
#    before it. Copying (rather than Worksheets.Add()) means the new
#    sheet starts out with the exact same cell styling (bold header /
#    index-column look) already in use.
# ---------------------------------------------------------------------
$totalSheetBefore = $wb.Worksheets.Item("总计")
$totalSheetBefore.Copy($totalSheetBefore)
$q1 = $wb.Worksheets.Item("总计 (2)")
$q1.Name = "2022-Q1"

# Header row: columns B..H.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    Set-TextValue $q1.Cells.Item(1, $col) $headers[$col - 2]
}
# Extend the existing header style (already on B1:D1) across the new E:H columns.
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

$q1Rows = @(
    @("159851", "华宝中证金融科技主题ETF", "3.16", "98.58", "4.04", "0.1277", 6),
    @("005062", "博时中证500指数增强A", "6.14", "90.01", "1.58", "0.0970", 5),
    @("516100", "华夏中证金融科技主题交易型开放式指数证券投资基金", "0.68", "96.91", "4.02", "0.0273", 6),
    @("005795", "博时中证500指数增强C", "1.14", "90.01", "1.58", "0.0180", 5),
    @("008112", "中泰中证500指数增强A", "0.61", "92.46", "1.41", "0.0086", 4),
    @("008113", "中泰中证500指数增强C", "0.46", "92.46", "1.41", "0.0065", 4)
)

$r = 2
foreach ($row in $q1Rows) {
    $q1.Cells.Item($r, 1).Value = $r - 2

    Set-TextValue $q1.Cells.Item($r, 2) $row[0]
    Set-TextValue $q1.Cells.Item($r, 3) $row[1]
    Set-TextValue $q1.Cells.Item($r, 4) $row[2]
    Set-TextValue $q1.Cells.Item($r, 5) $row[3]
    Set-TextValue $q1.Cells.Item($r, 6) $row[4]
    Set-TextValue $q1.Cells.Item($r, 7) $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]

    $r++
}
# Extend the existing index-column style (already on A2:A5) down across the
# two new rows (A6:A7).
$q1.Range("A5").Copy()
$q1.Range("A6:A7").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: push the existing data rows down by one and
#    write a new first row for 2022-Q1.
#    NOTE: must re-resolve the sheet by name now -- the older
#    `$totalSheetBefore` handle tracks the *position* that has since been
#    taken over by the freshly inserted "2022-Q1" sheet.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$existing = @(
    @("2021-Q4", 3, 0.08),
    @("2021-Q3", 2, 0.08),
    @("2021-Q2", 2, 0.09),
    @("2021-Q1", 4, 0.04)
)

for ($i = $existing.Length - 1; $i -ge 0; $i--) {
    $destRow = $i + 3
    $totalSheet.Cells.Item($destRow, 1).Value = $i + 1
    $totalSheet.Cells.Item($destRow, 2).Value = $existing[$i][0]
    $totalSheet.Cells.Item($destRow, 3).Value = $existing[$i][1]
    $totalSheet.Cells.Item($destRow, 4).Value = $existing[$i][2]
}

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 6
$totalSheet.Cells.Item(2, 4).Value = 0.29
